$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 34) appended to the stock/news-tone table.
$headlineText = 'spirit of halloweentown is calling festival ticket go on sale august passenger injured in icon of the sea water slide panel shatter kate spade outlet shopper say this crossbody bag is a timeless accessory and it on sale for only nvr inc announces share repurchase authorization car coffee and clothing are poised to get pricier with new tariff bitcoin news today bitcoin auction for dogwifhat meme cap hit k amid growing nft interest nvr inc announces share repurchase authorization top cryptos to invest before end of august best rakshabandhan tech gift under r taketwo q net booking up fran drescher wont run for sagaftra president a candidate line up for race gold miner earn double the profit of yet still undervalued trump warns of great depression a court weighs billion tariff rollback this small city ha managed to keep rent down since and it a quick train ride from nyc sec end lawsuit against ripple company to pay million fine reuters dow jones today stock surge a major index on track to post solid weekly gain apple extends rally gold at record high level investopedia investor alert morris kandinov llp investigating phh wai pcla and epwk stockholder are encouraged to contact the firm a florida woman posed a a nurse and treated thousand of unsuspecting patient official say gold price jump to record high following surprise u tariff ruling apple stock retake key level after dodging tariff bullet apple stock retake key level after dodging tariff bullet investor business daily trump pursues record billion settlement from ucla over antisemitism claim dwave quantum cash hoard fuel aggressive expansion analyst report parkerhannifin corp minsud announces option grant fat brand open first cobranded round table pizza and fatburger in california openstore shift focus from shopify aggregator to menswear brand run to samsung now for up to off a new tv score off your first order at small the purrfect cat meal service rosen a ranked and leading law firm encourages fiserv inc investor to secure counsel denver employee no longer with the city after incident at cannabis club rosen a ranked and leading law firm encourages fiserv inc investor to secure counsel before important deadline in security class action fi examining the future growgenerations earnings outlook a look at perasos upcoming earnings report a peek at teloss future earnings apple aapl confirms billion u investment after white house tease shore living welcome new realtor to team u global investor launch it smart beta sea etf on the mexican stock exchange joining jet and goau earnings preview for quest resource holding insight into quipt home medical upcoming earnings a preview of hf food group earnings examining the future virtras earnings outlook a look ahead gaia earnings forecast connex credit union data breach expose personal information murphy law firm investigates legal claim a glimpse of synchronoss technologiess earnings potential rosen a ranked and leading law firm encourages flywire corporation investor to secure uncovering potential microvast holdingss earnings preview analyst have this to say about elanco animal health an overview of aris water solutionss earnings deep dive into mirum pharmaceutical stock analyst perspective rating expert outlook abcellera biologics through the eye of analyst demystifying penn entertainment insight from analyst review what analyst rating have to say about iovance biotherapeutics tmobile u unusual option activity unpacking the latest option trading trend in cai unpacking the latest option trading trend in ford motor whats driving the market sentiment around affirm holding what doe the market think about paylocity holding pe ratio insight for jdcom taketwo interactive stock a deep dive into analyst perspective rating what doe the market think about ubs gr popular grocery store supplier file chapter bankruptcy how doe a k hardship withdrawal work and is it smart to take one ripple xrp tron trx poised for gain but trader are watching this crypto for a bigger run in samsung galaxy z fold what make it the best foldable phone cardanos ada lead in the race is in danger ruvi ai ruvi coinmarketcap listing get it industry spotlight a analyst call it the top summer play bofa reiterates buy on alphabet googl a gemini gain traction and kpis hold steady xrp news today xrp rally hit with institutional buying and legal victory in play rosen a ranked and leading law firm encourages flywire corporation investor to secure counsel before important deadline in security class action flyw fed president issue stern warning on lasting price increase trump clash with intels ceo is a crucial test for u chip making market update cmcsa cv expe msi pgr ph unh yum hims final deadline rosen a top ranked law firm encourages hims hers health inc fannie freddie stock offering could come a soon a this year barrons amazon amzn aws unveils amazon evsa new way to run vmware workload in the cloud recordhigh gold price push jeweler to offer gilded veneer to prevent sticker shock best ev penny stock to buy according to hedge fund royal expert reveal hypocrisy of william kates secret greek holiday electric vehicle supply equipment market size will attain usd billion by growing at cagr exclusive report by zion market research global electric vehicle supply equipment market size share trend analysis report electric vehicle supply equipment market size will attain usd billion by growing at cagr exclusive report by zion market research global electric vehicle supply equipment market size share trend analysis report san mateo burger popup start this weekend at coyote point best canadian gold stock to buy according to hedge fund doj asks court to release epstein grand jury evidencebut unclear how much may already be public trump approval rating drop in poll this week but rise slightly in retailer tariffbattered import volume to be weaker in gear for the great outdoors is up to off for this great sale at backcountry cool down anywhere save on this ultra cooling shark fan celebrate cbd day with off sitewide at sunday scaries qqq and friend hit high tech etf thrive despite trade turbulence parent trap star reunite in new movie freakier friday deutsche bank reiterates buy on microsoft msftai driving doubledigit efficiency gain ap business summarybrief at pm edt boar head plan to reopen troubled deli meat plant but report of sanitation problem persist boar head plan to reopen troubled deli meat plant but report of sanitation problem persist intel ceo fire back at trump demand to resign mcdonalds menu brings back classic character for shake meal nvidia nvda get overweight rating a ai capex stay strong into mexican authority accuse adidas of cultural appropriation over new sandal design gold future jump to record high a u imposes tariff on swiss gold bar import amd ai growth story gain steambarclays see room for more upside montgomery county democratic state rep call on senate republican to pas funding for septum friday assorted link opinion direct primary care provides medical treatment the oldfashioned way not sure whats more humiliating meidastouch pushing fake jd vance story or the obama bro who bought it report trump seek billion from ucla over antisemitism claim texas state rep isnt backing down from fighting overtly racist redistricting map investigation stalled while georgia state election board tried to change voting rule justice department subpoena new york ag james a it probe whether she violated trump right photo show israeli demonstrating for the release of hostage and ending the gaza war photo show israeli demonstrating for the release of hostage and ending the gaza war gov gavin newsom appoints new judge to bay area superior court kristi noem slam south park for petty and lazy spoof of her only the liberal nuclear energy now reactor on the moon and at denver airport bauer south carolina former lieutenant governor end gop primary challenge to sen graham bauer south carolina former lieutenant governor end gop primary challenge to sen graham bauer south carolina former lieutenant governor end gop primary challenge to sen graham watchdog group asks for probe of acquisition of qatari jet on trump behalf mandan approves preliminary budget with small property tax increase ap news summary at pm edt israel plan to take over gaza city stir fear for civilian and hostage the latest trump will host armenia and azerbaijan for a white house peace summit exclusive full video show leadup to brutal attack of elderly prolifers by man who dodged prison u global investor launch it smart beta sea etf on the mexican stock exchange exclusive federal law enforcement to begin interviewing unaccompanied migrant child in government custody texas republican vow to go nuclear in redistricting row trump action unlikely to affect fed independence trump tariff are killing f sale is there a better way to manage grizzly bear state urged to step up in matching federal momentum for developing nuclear energy trump considering use of military force against mexican drug cartel what doe occupying gaza accomplish for israel live now texas hold legislative session after democrat lawmaker depart doj asks court to release epstein grand jury evidencebut unclear how much may already be public trump approval rating drop in poll this week but rise slightly in plane cabin fill with toxic smoke a power bank catch fire over the ocean hour away from landing how closely do congressional delegation reflect how people vote not very how closely do congressional delegation reflect how people vote not very how closely do congressional delegation reflect how people vote not very how closely do congressional delegation reflect how people vote not very how closely do congressional delegation reflect how people vote not very how closely do congressional delegation reflect how people vote not very how closely do congressional delegation reflect how people vote not very how closely do congressional delegation reflect how people vote not very how closely do congressional delegation reflect how people vote not very president trump sign debanking executive order horizontal review for bank and credit union are next parent trap star reunite in new movie freakier friday boar head plan to reopen troubled deli meat plant but report of sanitation problem persist intel ceo fire back at trump demand to resign how closely do congressional delegation reflect how people vote not very how closely do congressional delegation reflect how people vote not very city can make do with fewer public employee florida updated agreement on handling detainee at alligator alcatraz but a month after it opened florida updated agreement on handling alligator alcatraz detainee but a month after it opened trump is asking the supreme court to bless stephen miller racial profiling white house pause public tour amid ballroom construction landrum police chief resigns after year of service minute before gunfire erupted at fort stewart the shooting suspect texted his family saying i love yall uncle say what to know about israel plan to retake gaza city what to know about israel plan to retake gaza city update labrynth launch the red tape index to help america build smarter faster other view millennials are avoiding divorce but their american dream ha changed house republican try to defend medicaid cutsand it go poorly dhs secretary kristi noem truth nuke a reporter pushing dem talking point disguised a question ice family separation are forcing child to parent themselves ap news summary at pm edt columnist helen ubinas my philadelphia inquirer farewell lt governor kounalakis announces switch for governor race to treasurer thomas financial group and phoenix lender service close million usda business possible outcome of a gerrymandering battle royale reason republican redistricting power grab might backfire doe your state have official drink if so what is it here the answer warren call on fed to activate additional capital charge special prosecutor ag bondi assigns to letitia james and adam schiff should scare the hell out of them trumpappointed judge toss contempt finding against administration doj informs judge it will ask supreme court to quickly rule on trump birthright citizenship order epa terminates billion solar for all program what to know about past meeting between putin and his american counterpart germany halt military export to israel for use in gaza amid outcry over netanyahu plan when real life feel like a jurassic park moment pete davidson thought the starstudded audience at the snl special wa terrible just peachy fema tell staff to work for ice or risk getting fired letter to the editor leave the judging to jesus tarantula overcome southwest during mating season right now the worst thing you can be is a racist marc maron fire back at howie mandels comment using donald trump a an example dam if you do dam if dont private berk lake murky future opinion opinion to keep finest city title san diego must clean up it act russia mediterranean flotilla is in a sorry state aquaculture can help produce more u seafood pennsylvania advocate call for shutdown of moshannon after ice detainee death putin say hell end war with ukraine with one major concession trump double reward to million for arrest of venezuela president to face u drug charge blue state throw tantrum over trump food stamp fraud crackdown couple welcome baby born from yearold embryo sen kelly bill help military medic enter civilian career b bombed a smiley face into the ground because it ran out of target jd vance go fishing for trout with uk lammy germany halt military export to israel for use in gaza amid outcry over netanyahu plan germany halt military export to israel for use in gaza amid outcry over netanyahu plan appeal court remove boasbergs contempt order in trump admin deportation case why trump tariff hit poor country with higher rate than rich one diana taurasi endorses caitlin clark joining shortlist of wnba player with signature shoe mixedbreed dog buddy is a sweetheart ready for a welcoming new home a top federal reserve official say dour job data back the case for rate cut a top federal reserve official say dour job data back the case for rate cut a top federal reserve official say dour job data back the case for rate cut common lead a strong opening night for san jose jazz summer fest a top federal reserve official say dour job data back the case for rate cut ap news morton meyerson texas businessman and perot lieutenant dy at why xiaomi succeeded at making electric car but apple failed jim cramer delivers straight talk on tricky sp market top crypto token under that could see big gain in fastestgrowing cryptos in blockdag toncoin cronos tron gaining fast top cryptos under to buy nowozak ai lead the list while nearing m presale momentum and x token growth potential cardano price prediction ada eye x growth but this cheaper crypto might soar x in why clean energy fuel stock raced nearly higher today here why joby aviation stock flew higher in july china ev tipping point racing from to new sale in record time rethinking state government presence in baltimore editorial massachusetts makarevich understandable separation agreement language aid employer in unpaid wage case building growth lane that work in any market weekend binge best family series to watch on hotstar from low market cap to big ambition most promising cryptos in landlord say ruby lius plan for bay property defies commercial common sense rosie odonnell claim abc gonna cancel the view michigan governor tell president trump tariff are hurting auto job fda official return to agency after loomerled ouster cnn lula and putin discus peace in ukraine before u summit after the mustang ford is teasing the return of another beloved car a an ev gizmodo apple ai momentum is building here what it mean for investor top crypto gainer today blockdag stellar uniswap hedera are the one to watch ukraine drone revolution ha exposed a u military vulnerability how to trade in corporate life for your next big idea best meme coin to buy new pepe coin rival is growing faster than pepe when it gained billion market cap in day court rule pension fund position wa not so baseless a to mandate an award of attorney fee new art cafe community space to open in camp washington this summer whistleblower tie clinton campaign to fake russia hack sperry trump meeting with putin pose high stake for the ukraine war and his legacy investor alert pomerantz law firm investigates claim on behalf of investor of agilon health inc agl investor alert pomerantz law firm investigates claim on behalf of investor of inspire medical system inc insp dol resurrects paid program employer can selfreport resolve violation rosen recognized investor counsel encourages lockheed martin corporation investor to secure investor alert pomerantz law firm investigates claim on behalf of investor of national grid rosen recognized investor counsel encourages lockheed martin corporation investor to secure counsel before important deadline in security class action lmt trump pick for the fed fuel an existential threat a central bank independence is targeted jpmorgan say fortune what to know about instagram map location sharing and privacy concern the new york time a gpt launch openai loyalist pushed altman for gptos return citing performance and reliability flipkart freedom sale huge discount on ipad oneplus pad samsung tab top presale cryptos in this month that are breaking funding record cold wallet dexboss aurealone lightchain ai taxloss harvesting part iii investment strategy u fed bowman latest job data stiffens support for three rate cut in california rooftop solar rule in limbo after state supreme court ruling top altcoins to buy blockchainfx presale at might dwarf polkadot and cosmos gain speech by vice chair for supervision bowman on the economic outlook and community banking federal reserve board gov irs white house clashed over immigrant data before tax chief wa ousted u fed bowman latest job data stiffens support for three rate cut in reuters fed bowman call for community bank reform favor rate cut these are dead decision zelensky reject ceding territory after trump talk landswap fed bowman back a september rate cut the wall street journal thumzup medium corporation announces update to the term of it proposed public offering im what is the perfect growth etf to pair with voo for longterm wealth reddit stock pop after a blowout quarter but is it a buy best crypto to buy this month magacoin finance gain whale support from solana and doge polk county chamber of commerce ha new home from hollywood to silicon valley how eventbrites ceo bet k on a phone closet startup is qualcomm the best semiconductor stock to buy right now best crypto to buy now blockdag shiba inu cardano xrp stand out fpl reach potential rate settlement deep dive american havent saved for retirement state are creating automatic saving plan ai can solve the genz and millennial engagement crisis put yourself in front of the next major investing trend trump tariff are making international stock market great again reason why ozak ai could be the best investment this august x growth still on the table investor alert pomerantz law firm investigates claim on behalf of investor of barnes noble education inc bned why burger cost so much right now trump announces white house press conference to stop violent crime in washington dcdespite historically low rate whitmer told trump in private that michigan auto job depend on a tariff change of course whitmer told trump in private that michigan auto job depend on a tariff change of course how a cia hit on al qaeda ensnared a u citizen in afghanistan byd ev struck by lightning hold up fine your daily dose of insightful news and perspective investor alert pomerantz law firm investigates claim on behalf of investor of exelixis inc exel ai can be used to develop biological and chemical weapon investor alert pomerantz law firm investigates claim on behalf of investor of anika investor alert pomerantz law firm investigates claim on behalf of investor of anika therapeutic inc anik investor alert pomerantz law firm investigates claim on behalf of investor of neogenomics inc neo investor alert pomerantz law firm investigates claim on behalf of investor of praxis precision medicine inc prax how compounding return can help you retire a millionaire even on a modest income investor alert pomerantz law firm investigates claim on behalf of investor of albany international corp ain top new cryptocurrencies to buy in top presale coin to buy today before price jump in the next stage investor alert pomerantz law firm investigates claim on behalf of investor of tronox holding plc trox investor alert pomerantz law firm investigates claim on behalf of investor of confluent inc cflt crew splash down in pacific after return from iss if youd invested in nvidia stock year ago here how much youd have today vinay prasad return to the fda week after his ouster statnewscom farmer consider the good and bad of incorporating artificial intelligence on the farm farmer consider the good and bad of incorporating artificial intelligence on the farm only a blink left to buy in arctic pablo shock market with b toke burned a shiba inu and pudgy penguin rally bathroom tool that save money and water all you need is this simple upgrade reader view essentia aspirus need to do more for duluth laguna niguel man pleads guilty to firing semiautomatic rifle at sheriff helicopter justice department target new york attorney general a trump foe here what to know justice department target new york attorney general a trump foe here what to know man suckerpunched during argument at eureka sport bar charge filed end of racial consent decree poised to change federal hiring sin city economic slump signal wider slowdown the soviet union yakk fighter wa too powerful for it own design wnba coach say female athlete have been subjected to adult toy projectile for century wut trump team push for ouster of top iea official california poverty crisis will hamper newsoms national climb avoid buying this canned coconut milk brand our taste test show it the absolute worst security flaw in web server enables remote code execution cve london police arrest people a propalestinian protester defy new law london police arrest people a propalestinian protester defy new law trump narcos strategy collides with mexico strong opposition to u deployment oc man sentenced to year for la street takeover that left nursing student dead extrump lawyer call president middecade redistricting scheme more brazen than watergate welcome to the next phase of american heist secretary duffy host press conference on one big beautiful bill act and air traffic control trump warns of new great depression if court block his tariff plan it would be all over again why did judge throw out conviction in murder of chris paul grandfather in winstonsalem beto orourke trash texas republican these motherfkers are panicking beto orourke trash texas republican these motherfkers are panicking trump to host news conference to stop violent crime and make dc one of the safest city in the world mexican president rule out trump reported military plan against mexico drug cartel kelly army arrival worried the town fda official return to agency after loomerled ouster illinois gov jb pritzker forgets map exist when slamming antidemocracy gopers over redistricting nyc limit housing discrimination based on criminal background is criminal history history breastfeeding at work redefined puerto rico new code usher in major change to save animal county must mobilize more volunteer to save animal county must mobilize more volunteer diddy lawyer mark geragos say diddy would be trump supporter if pardoned clip of the week old man wanders around on a roofand more the soviet union insane first attempt at an aircraft carrier ag paxton call on california to aid texas in arresting house democrat that left the state ai the road to utopia or dystopia man charged with raping unconscious woman multiple time in irvine trump meeting with putin pose high stake for the ukraine war and his legacy homeland security secretary kristi noem blast illinois leader during visit ambassador mike huckabee slam british pm keir starmer for for condemning israel did uk surrender to nazi and drop food to them doj seek to unseal jeffrey epstein ghislaine maxwell grand jury record chairlift collapse at a resort and leaf at least hurt including person who sustained a spinal cord injury watt rebellion th anniversary renews call for kerner report plan to fight poverty court order reinstatement and back pay for pittsburgh officer fired after taser arrest death restored nagasaki bell ring in year since abomb strong march demand britain break tie with israel howie carr suffolk high sheriff steven tompkins earns his democrat credential government house repair wont break taxpayer reader commentary california focus require id from federal agent end some fear israel face growing global condemnation over military expansion in gaza founder of trump burger chain could be deported after ice arrest bernie sander kick off fighting the oligarchy tour in west virginia the dred scott dissent lincoln loved the three strategic bomber the u air force cant live without la county evacuation warning lifted a canyon fire continues burning national perspective newton third law rule washington politics artifact in u museum explain the outsized role of racism in the nation history did i miss a memo adam schiffs lawyer say mortgage fraud allegation have been debunked trump announces white house press conference to stop violent crime in washington dcdespite historically low rate fda leader resigned after push from big pharma now he back live now texas house hold legislative session on aug for democrat lawmaker to come back investor alert pomerantz law firm investigates claim on behalf of investor of exelixis inc exel hundred defy police brutality to oppose palestine action ban investor alert pomerantz law firm investigates claim on behalf of investor of anika therapeutic inc anik investor alert pomerantz law firm investigates claim on behalf of investor of neogenomics inc neo investor alert pomerantz law firm investigates claim on behalf of investor of praxis precision medicine inc prax the nation capital wait for trump next move a a federal takeover threat loom investor alert pomerantz law firm investigates claim on behalf of investor of albany international corp ain trump planned computer chip tariff spark confusion among business and trading partner investor alert pomerantz law firm investigates claim on behalf of investor of tronox holding plc trox investor alert pomerantz law firm investigates claim on behalf of investor of confluent inc cflt how to keep hezbollah away from power gop seek to game system thats not a noble cause reader commentary larimer county restaurant inspection made aug iran wont allow trumpbacked azerbaijan corridor top aide say he wa on the fbi most wanted list now he running for congress this physicianscientist is taking on trump on behalf of disadvantaged community we are at war bring it on democrat ready to fight dirty to stop trump gerrymandering is despicable chase empsall the time is now to speak out rick hausman the f lightning ii is still the greatest fighter jet in the world france sends crew to keep fatal wildfire from reigniting shaw walz shouldnt run for reelection and fargo poor communication david cone drop new country album if you dont count the money im the richest man in town grand haven twp to discus pair of development past page for august to accused in death of manchester emt retired firefighter ordered held on unrelated charge susan shelley good riddance to cafe penalty lawmaker urge meta to shut down instagram map abysmal at protecting child politics without shame gerrymandering make hypocrisy a political punch line judge block beto orourke from financially supporting texas democrat who left the state local view like caesar rome u feeling more vulnerable pomerantz law firm announces the filing of a class action against tesla inc and certain officer tsla bonior a tribute to the courageous ukrainian people london police arrest demonstrator supporting the banned palestine action group pomerantz law firm announces the filing of a class action against irobot corporation and certain officer irbt brickner navalny an american hero letter to the editor geothermal is way to go letter to the editor pritzker just a hypocrite white house mull inviting zelensky to trumpputin talk in alaska d printing and ai cut nuclear reactor build time to day google test ai finance platform with chatbot and live insight stanford forgoes m cal grant to keep legacy admission boar head to reopen virginia plant after fatal listeria outbreak concern why cold wallet top the list of top crypto to buy with x roi outpacing shiba inus and toncoins gain twenty suspected mossad spy arrested in iran top meme coin in august with massive roi potential aveo capital partner llc decrease position in exxon mobil corporation nysexom shaquille oneal admits family mistake left him lonely in squarefoot home i know i messed up and when i didnt have that i wa lost owner connected to monster energy announce plan to expand thrifty ice cream gasket market is forecasted to reach u billion in say stratview research microelectronics vacuum valve market is forecasted to reach u billion in say stratview research semiconductor fluid conveyance market is forecasted to reach u billion in say stratview research semiconductor seal market is forecasted to reach u billion in say stratview research the debt and deficit problem isnt what you think microelectronics vacuum valve market is forecasted to reach u billion in say stratview research forget shiba inu shib and pepe coin pepe the best meme coin to buy in is still very undervalued semiconductor vacuum valve market is forecasted to reach u billion in say stratview research semiconductor fluid conveyance market is forecasted to reach u billion in say stratview research semiconductor seal market is forecasted to reach u billion in say stratview research gasket market is forecasted to reach u billion in say stratview research last hour to grab blockdag at before price change xrp hype reach key level semiconductor vacuum valve market is forecasted to reach u billion in say stratview research weapon horror film score a box office victory weapon horror film score a box office victory tesla cybertruck is suddenly sold out sort of gizmodo cobblestone capital advisor llc ny decrease stock position in oracle corporation nyseorcl ccm investment group llc purchase new position in chevron corporation nysecvx third ave outlook brightening cbre say largest u retailer to pay million for overcharging customer newsnation new offshore wind study confirms the obvious a malevolently incompetent clown is steering the u clown car financial literacy key term every firstgeneration american should know luxury houseboat rental a growing trend in virginia beach solana price prediction suggests modest gainsozak ai could leave it behind strengthening family community llc acquires share of mcdonalds corporation nysemcd accenture plc nyseacn share bought by covington investment advisor inc coin gaining steam the best altcoin to buy now for major growth largest stock exchange in the world by market cap in designing a weekly financial health check for your truck or fleet best ethereum ecosystem coin by market cap in connor clark lunn investment management ltd acquires new share in eaton corporation plc nyseetn eye on atlantic basin a tropical development likely next week more than half of industry are already shedding worker a telling sign thats accompanied past recession top economist say fortune united flight ua is declaring an emergency and diverting to new york airlive chevron corporation nysecvx share acquired by first american trust fsb linscomb wealth inc ha million stock holding in visa inc nysev northeast financial group inc ha holding in eli lilly and company nyselly business people coborns exec david best tapped to lead cub food nina roza longquan the dragon spring and dark chocolate scoop prize at locarno pro industry award king luther capital management corp trim stake in exxon mobil corporation nysexom a robot cop named parker spark controversy in wealthy maryland county the washington post inflation fed leadership and other key thing to watch this week blockdag achieves m presale milestone in record time solana hold at bnb near bessent say tariff should melt if trade rebalances nik'

$ws.Cells.Item(34, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 1).Value = 45880
$ws.Cells.Item(34, 2).Value = 6373.4501953125
$ws.Cells.Item(34, 3).Value = 6407.25
$ws.Cells.Item(34, 4).Value = 6364.06005859375
$ws.Cells.Item(34, 5).Value = 6389.669921875
$ws.Cells.Item(34, 6).Value = 4652400000
$ws.Cells.Item(34, 7).Value = -0.0025041278217863
$ws.Cells.Item(34, 8).Value = -1
$ws.Cells.Item(34, 9).Value = $headlineText
$ws.Cells.Item(34, 10).Value = $headlineText
$ws.Cells.Item(34, 11).Value = 0.149
$ws.Cells.Item(34, 12).Value = 0.12
$ws.Cells.Item(34, 13).Value = 0.73
$ws.Cells.Item(34, 14).Value = 0.9993
$ws.Cells.Item(34, 15).Value = 208
$ws.Cells.Item(34, 16).Value = 179
$ws.Cells.Item(34, 17).Value = 0
$ws.Cells.Item(34, 18).Value = 344
$ws.Cells.Item(34, 19).Value = 130
$ws.Cells.Item(34, 20).Value = 424
$ws.Cells.Item(34, 21).Value = 275
$ws.Cells.Item(34, 22).Value = 114
$ws.Cells.Item(34, 23).Value = 64
$ws.Cells.Item(34, 24).Value = 127

Write-Host "Row 34 written"
